$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("A2").Value = "40e526d7-263a-4f74-b935-1359b190b926"
$ws.Range("B2").Value = "2018-09-10 04:59:51.285837"
$ws.Range("C2").Value = "2018-09-10 05:00:09.798837"
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 26
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = "Casting"
$ws.Range("H2").Value = "Casting: OK"

# Row 3
$ws.Range("A3").Value = "40e526d7-263a-4f74-b935-1359b190b926"
$ws.Range("B3").Value = "2018-09-10 05:06:23.535176"
$ws.Range("C3").Value = "2018-09-10 05:06:59.399146"
$ws.Range("D3").Value = 28
$ws.Range("E3").Value = 25.8
$ws.Range("F3").Value = 23
$ws.Range("G3").Value = "Casting"
$ws.Range("H3").Value = "Casting: Bad"

# Row 4
$ws.Range("A4").Value = "40e526d7-263a-4f74-b935-1359b190b926"
$ws.Range("B4").Value = "2018-09-10 05:14:52.426060"
$ws.Range("C4").Value = "2018-09-10 05:15:53.965341"
$ws.Range("D4").Value = 27.4
$ws.Range("E4").Value = 24.8
$ws.Range("F4").Value = 24.6
$ws.Range("G4").Value = "Casting"
$ws.Range("H4").Value = "Casting: Bad"

# Row 5
$ws.Range("A5").Value = "Test"
$ws.Range("B5").Value = "2018-09-10 05:32:00.874461"
$ws.Range("C5").Value = "2018-09-10 05:32:25.816646"
$ws.Range("D5").Value = 26.4
$ws.Range("E5").Value = 23.8
$ws.Range("F5").Value = 26.2
$ws.Range("G5").Value = "Casting"
$ws.Range("H5").Value = "Casting: Bad"

# Row 6 - D,E,F are stored as strings (text), not numbers
$ws.Range("A6").Value = "Test"
$ws.Range("B6").Value = "2018-09-10 05:32:00.874461"
$ws.Range("C6").Value = "2018-09-10 05:33:02.706576"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "26.2"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "25.8"
$ws.Range("E6").Style = "Normal"

$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "24.8"
$ws.Range("F6").Style = "Normal"

$ws.Range("G6").Value = "Casting"
$ws.Range("H6").Value = "Casting: Bad"

# Reset selection to A1 (matches the committed worksheet view)
$ws.Range("A1").Select() | Out-Null
